$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 294 (shifts existing rows 294:369 down to 295:370)
$ws.Rows("294:294").Insert()

# Populate the new row 294 with the new record's data.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R are constant across all records in this sheet.
$ws.Range("A294").Value = 3
$ws.Range("B294").Value = "Femacal de La Calera"
$ws.Range("C294").Value = "Coquimbo"
$ws.Range("D294").Value = 44722
$ws.Range("E294").Value = 5
$ws.Range("F294").Value = 100112040
$ws.Range("G294").Value = "Cilantro"
$ws.Range("H294").Value = "Sin especificar"
$ws.Range("I294").Value = "Primera"
$ws.Range("J294").Value = 310
$ws.Range("K294").Value = 3000
$ws.Range("L294").Value = 3500
$ws.Range("M294").Value = 3258
$ws.Range("N294").Value = "$/docena de atados (3 kilos)"
$ws.Range("O294").Value = "Provincia de Quillota"
$ws.Range("P294").Value = 1086
$ws.Range("Q294").Value = 3
$ws.Range("R294").Value = "Hortaliza"
